# Applies the "dashboard sketch" edit:
#  - Rectangle 3 (ABV input box): add "Number of drinks" and "Time drinking"
#    paragraphs right after the existing "Hours drinking" paragraph.
#  - Rectangle 6 / Rectangle 7 ("Disclaimer?" boxes): append
#    " (include the BRAD information)" to the text.
#  - Rectangle 8 (Map box): change "Other?" to "Bars".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    throw "Shape not found: $name"
}

function Get-ParagraphIndex($textRange, $expectedText) {
    $count = $textRange.Paragraphs().Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $textRange.Paragraphs($i, 1)
        if ($para.Text.TrimEnd("`r") -eq $expectedText) {
            return $i
        }
    }
    throw "Paragraph not found: $expectedText"
}

# --- Rectangle 3: insert two new paragraphs after "Hours drinking" ---
$inputBox = Get-ShapeByName $s "Rectangle 3"
$inputTr = $inputBox.TextFrame.TextRange
$hoursIdx = Get-ParagraphIndex $inputTr "Hours drinking"
$hoursPara = $inputTr.Paragraphs($hoursIdx, 1)
$hoursPara.InsertAfter("`rNumber of drinks`rTime drinking") | Out-Null

# --- Rectangle 6: "Disclaimer?" -> "Disclaimer? (include the BRAD information)" ---
$disclaimer1 = Get-ShapeByName $s "Rectangle 6"
$disclaimer1.TextFrame.TextRange.Text = "Disclaimer? (include the BRAD information)"

# --- Rectangle 7: "Disclaimer?" -> "Disclaimer? (include the BRAD information)" ---
$disclaimer2 = Get-ShapeByName $s "Rectangle 7"
$disclaimer2.TextFrame.TextRange.Text = "Disclaimer? (include the BRAD information)"

# --- Rectangle 8 (Map): "Other?" -> "Bars" ---
$mapBox = Get-ShapeByName $s "Rectangle 8"
$mapTr = $mapBox.TextFrame.TextRange
$otherIdx = Get-ParagraphIndex $mapTr "Other?"
$otherPara = $mapTr.Paragraphs($otherIdx, 1)
$otherPara.Text = "Bars"
